$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the input pulse measurement values
$ws.Range("D3").Value = 1000
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 2000
$ws.Range("D8").Value = 24
$ws.Range("D13").Value = 2
$ws.Range("D20").Value = 100000
$ws.Range("D21").Value = 27000

# Update the selected cell on the sheet
$ws.Range("D9").Select()
